$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2265, 1).Value = 'Buying Opportunity'
$ws.Cells.Item(2265, 2).Value = 'support Zone'
$ws.Cells.Item(2265, 3).Value = 'long buildup'
$ws.Cells.Item(2265, 4).Value = 'Short buildup'
$ws.Cells.Item(2265, 5).Value = 'FII ENTERING'

$ws.Cells.Item(2266, 1).Value = 'AVROIND'
$ws.Cells.Item(2266, 2).Value = 'APTUS'
$ws.Cells.Item(2266, 4).Value = 'AARTIIND'
$ws.Cells.Item(2266, 6).Value = 125.04
$ws.Cells.Item(2266, 7).Value = 305.35
$ws.Cells.Item(2266, 9).Value = 621.15

$ws.Cells.Item(2267, 1).Value = 'BALKRISHNA'
$ws.Cells.Item(2267, 2).Value = 'AVALON'
$ws.Cells.Item(2267, 4).Value = 'ATUL'
$ws.Cells.Item(2267, 6).Value = 24.13
$ws.Cells.Item(2267, 7).Value = 459.75
$ws.Cells.Item(2267, 9).Value = 7635

$ws.Cells.Item(2268, 1).Value = 'BALRAMCHIN'
$ws.Cells.Item(2268, 2).Value = 'AVTNPL'
$ws.Cells.Item(2268, 4).Value = 'BATAINDIA'
$ws.Cells.Item(2268, 6).Value = 511.4
$ws.Cells.Item(2268, 7).Value = 85.23
$ws.Cells.Item(2268, 9).Value = 1402.8

$ws.Cells.Item(2269, 1).Value = 'CMSINFO'
$ws.Cells.Item(2269, 2).Value = 'BAYERCROP'
$ws.Cells.Item(2269, 4).Value = 'BHEL'
$ws.Cells.Item(2269, 6).Value = 556.5
$ws.Cells.Item(2269, 7).Value = 6214.05
$ws.Cells.Item(2269, 9).Value = 289.95

$ws.Cells.Item(2270, 1).Value = 'DALMIASUG'
$ws.Cells.Item(2270, 2).Value = 'BLS'
$ws.Cells.Item(2270, 4).Value = 'IDEA'
$ws.Cells.Item(2270, 6).Value = 397.65
$ws.Cells.Item(2270, 7).Value = 370.95
$ws.Cells.Item(2270, 9).Value = 15.47

$ws.Cells.Item(2271, 1).Value = 'DEEPENR'
$ws.Cells.Item(2271, 2).Value = 'CARERATING'
$ws.Cells.Item(2271, 4).Value = 'IRCTC'
$ws.Cells.Item(2271, 6).Value = 183.26
$ws.Cells.Item(2271, 7).Value = 971
$ws.Cells.Item(2271, 9).Value = 918.45

$ws.Cells.Item(2272, 1).Value = 'DHAMPURSUG'
$ws.Cells.Item(2272, 2).Value = 'CHEVIOT'
$ws.Cells.Item(2272, 4).Value = 'LT'
$ws.Cells.Item(2272, 6).Value = 202.73
$ws.Cells.Item(2272, 7).Value = 1406.2
$ws.Cells.Item(2272, 9).Value = 3551.8

$ws.Cells.Item(2273, 1).Value = 'EMAMILTD'
$ws.Cells.Item(2273, 2).Value = 'CONFIPET'
$ws.Cells.Item(2273, 4).Value = 'RECLTD'
$ws.Cells.Item(2273, 6).Value = 816.9
$ws.Cells.Item(2273, 7).Value = 81.37
$ws.Cells.Item(2273, 9).Value = 568.95

$ws.Cells.Item(2274, 1).Value = 'FORTIS'
$ws.Cells.Item(2274, 2).Value = 'CONTROLPR'
$ws.Cells.Item(2274, 6).Value = 517.1
$ws.Cells.Item(2274, 7).Value = 802.05

$ws.Cells.Item(2275, 1).Value = 'GULFOILLUB'
$ws.Cells.Item(2275, 2).Value = 'DCMSHRIRAM'
$ws.Cells.Item(2275, 6).Value = 1336.1
$ws.Cells.Item(2275, 7).Value = 1115.4

$ws.Cells.Item(2276, 1).Value = 'HEG'
$ws.Cells.Item(2276, 2).Value = 'GEOJITFSL'
$ws.Cells.Item(2276, 6).Value = 2169.75
$ws.Cells.Item(2276, 7).Value = 101.13

$ws.Cells.Item(2277, 1).Value = 'ISEC'
$ws.Cells.Item(2277, 2).Value = 'GRSE'
$ws.Cells.Item(2277, 6).Value = 785.9
$ws.Cells.Item(2277, 7).Value = 1951.3

$ws.Cells.Item(2278, 1).Value = 'KAYNES'
$ws.Cells.Item(2278, 2).Value = 'GSLSU'
$ws.Cells.Item(2278, 6).Value = 4702
$ws.Cells.Item(2278, 7).Value = 182.71

$ws.Cells.Item(2279, 1).Value = 'LOTUSEYE'
$ws.Cells.Item(2279, 2).Value = 'IMPAL'
$ws.Cells.Item(2279, 6).Value = 77.45
$ws.Cells.Item(2279, 7).Value = 1247.9

$ws.Cells.Item(2280, 1).Value = 'MADHAV'
$ws.Cells.Item(2280, 2).Value = 'KABRAEXTRU'
$ws.Cells.Item(2280, 6).Value = 52.39
$ws.Cells.Item(2280, 7).Value = 383.05

$ws.Cells.Item(2281, 1).Value = 'MAGADSUGAR'
$ws.Cells.Item(2281, 2).Value = 'KRITI'
$ws.Cells.Item(2281, 6).Value = 752.65
$ws.Cells.Item(2281, 7).Value = 227.57

$ws.Cells.Item(2282, 1).Value = 'MAXESTATES'
$ws.Cells.Item(2282, 2).Value = 'LAMBODHARA'
$ws.Cells.Item(2282, 6).Value = 589.85
$ws.Cells.Item(2282, 7).Value = 174.9

$ws.Cells.Item(2283, 1).Value = 'PLAZACABLE'
$ws.Cells.Item(2283, 2).Value = 'LAURUSLABS'
$ws.Cells.Item(2283, 6).Value = 90.28
$ws.Cells.Item(2283, 7).Value = 419.85

$ws.Cells.Item(2284, 1).Value = 'PRAJIND'
$ws.Cells.Item(2284, 2).Value = 'MAZDA'
$ws.Cells.Item(2284, 6).Value = 736.4
$ws.Cells.Item(2284, 7).Value = 1246.85

$ws.Cells.Item(2285, 1).Value = 'PRESTIGE'
$ws.Cells.Item(2285, 2).Value = 'MFSL'
$ws.Cells.Item(2285, 6).Value = 1778.05
$ws.Cells.Item(2285, 7).Value = 1052.55

$ws.Cells.Item(2286, 1).Value = 'RANASUG'
$ws.Cells.Item(2286, 2).Value = 'MSTCLTD'
$ws.Cells.Item(2286, 6).Value = 24.01
$ws.Cells.Item(2286, 7).Value = 851.05

$ws.Cells.Item(2287, 1).Value = 'RANEHOLDIN'
$ws.Cells.Item(2287, 2).Value = 'PCBL'
$ws.Cells.Item(2287, 6).Value = 1723.1
$ws.Cells.Item(2287, 7).Value = 374.95

$ws.Cells.Item(2288, 2).Value = 'ROHLTD'
$ws.Cells.Item(2288, 7).Value = 344.6

$ws.Cells.Item(2289, 2).Value = 'RRKABEL'
$ws.Cells.Item(2289, 7).Value = 1590.35

$ws.Cells.Item(2290, 2).Value = 'SADHNANIQ'
$ws.Cells.Item(2290, 7).Value = 70.61

$ws.Cells.Item(2291, 1).Value = '13/08/2024'

$ws.Cells.Item(2292, 1).Value = 'Buying Opportunity'
$ws.Cells.Item(2292, 2).Value = 'support Zone'
$ws.Cells.Item(2292, 3).Value = 'long buildup'
$ws.Cells.Item(2292, 4).Value = 'Short buildup'
$ws.Cells.Item(2292, 5).Value = 'FII ENTERING'

$ws.Cells.Item(2293, 1).Value = 'FMNL'
$ws.Cells.Item(2293, 2).Value = 'ACC'
$ws.Cells.Item(2293, 4).Value = 'GMRINFRA'
$ws.Cells.Item(2293, 6).Value = 6.55
$ws.Cells.Item(2293, 7).Value = 2281.95
$ws.Cells.Item(2293, 9).Value = 92.73

$ws.Cells.Item(2294, 1).Value = 'ISEC'
$ws.Cells.Item(2294, 2).Value = 'ACL'
$ws.Cells.Item(2294, 4).Value = 'HAL'
$ws.Cells.Item(2294, 6).Value = 805.65
$ws.Cells.Item(2294, 7).Value = 92.13
$ws.Cells.Item(2294, 9).Value = 4661.7

$ws.Cells.Item(2295, 1).Value = 'PRESTIGE'
$ws.Cells.Item(2295, 2).Value = 'ADANIENT'
$ws.Cells.Item(2295, 4).Value = 'MUTHOOTFIN'
$ws.Cells.Item(2295, 6).Value = 1787.5
$ws.Cells.Item(2295, 7).Value = 3040.1
$ws.Cells.Item(2295, 9).Value = 1816.45

$ws.Cells.Item(2296, 2).Value = 'AKZOINDIA'
$ws.Cells.Item(2296, 7).Value = 3194.85

$ws.Cells.Item(2297, 2).Value = 'ALOKINDS'
$ws.Cells.Item(2297, 7).Value = 25.05

$ws.Cells.Item(2298, 2).Value = 'AMNPLST'
$ws.Cells.Item(2298, 7).Value = 314.65

$ws.Cells.Item(2299, 2).Value = 'APOLLO'
$ws.Cells.Item(2299, 7).Value = 106.69

$ws.Cells.Item(2300, 2).Value = 'AWL'
$ws.Cells.Item(2300, 7).Value = 352.3

$ws.Cells.Item(2301, 2).Value = 'AXITA'
$ws.Cells.Item(2301, 7).Value = 24.17

$ws.Cells.Item(2302, 2).Value = 'BANARISUG'
$ws.Cells.Item(2302, 7).Value = 3067.2

$ws.Cells.Item(2303, 2).Value = 'BDL'
$ws.Cells.Item(2303, 7).Value = 1317.8

$ws.Cells.Item(2304, 2).Value = 'BHARATFORG'
$ws.Cells.Item(2304, 7).Value = 1567.2

$ws.Cells.Item(2305, 2).Value = 'CENTUM'
$ws.Cells.Item(2305, 7).Value = 1509.05

$ws.Cells.Item(2306, 2).Value = 'CONCOR'
$ws.Cells.Item(2306, 7).Value = 951.4

$ws.Cells.Item(2307, 2).Value = 'CONFIPET'
$ws.Cells.Item(2307, 7).Value = 80.94

$ws.Cells.Item(2308, 2).Value = 'COSMOFIRST'
$ws.Cells.Item(2308, 7).Value = 874.15

$ws.Cells.Item(2309, 2).Value = 'GEOJITFSL'
$ws.Cells.Item(2309, 7).Value = 100.3

$ws.Cells.Item(2310, 2).Value = 'GILLANDERS'
$ws.Cells.Item(2310, 7).Value = 86.04000000000001

$ws.Cells.Item(2311, 2).Value = 'GLOBAL'
$ws.Cells.Item(2311, 7).Value = 171.53

$ws.Cells.Item(2312, 2).Value = 'GODREJAGRO'
$ws.Cells.Item(2312, 7).Value = 795.5

$ws.Cells.Item(2313, 2).Value = 'GRASIM'
$ws.Cells.Item(2313, 7).Value = 2512.4

$ws.Cells.Item(2314, 2).Value = 'HDFCLOWVOL'
$ws.Cells.Item(2314, 7).Value = 19.95

$ws.Cells.Item(2315, 2).Value = 'HINDWAREAP'
$ws.Cells.Item(2315, 7).Value = 372.35

$ws.Cells.Item(2316, 2).Value = 'HLVLTD'
$ws.Cells.Item(2316, 7).Value = 19.02

$ws.Cells.Item(2317, 2).Value = 'HPAL'
$ws.Cells.Item(2317, 7).Value = 95.31

$ws.Cells.Item(2318, 2).Value = 'ICIL'
$ws.Cells.Item(2318, 7).Value = 361.35

$ws.Cells.Item(2319, 2).Value = 'IMPAL'
$ws.Cells.Item(2319, 7).Value = 1229

$ws.Cells.Item(2320, 2).Value = 'INFRABEES'
$ws.Cells.Item(2320, 7).Value = 935.72

$ws.Cells.Item(2321, 2).Value = 'J&KBANK'
$ws.Cells.Item(2321, 7).Value = 108.02

$ws.Cells.Item(2322, 2).Value = 'JAGRAN'
$ws.Cells.Item(2322, 7).Value = 89.11

$ws.Cells.Item(2323, 2).Value = 'JAYSREETEA'
$ws.Cells.Item(2323, 7).Value = 122.53

$ws.Cells.Item(2324, 2).Value = 'JUBLINGREA'
$ws.Cells.Item(2324, 7).Value = 654.45

$ws.Cells.Item(2325, 2).Value = 'KOHINOOR'
$ws.Cells.Item(2325, 7).Value = 40.08

$ws.Cells.Item(2326, 2).Value = 'LICI'
$ws.Cells.Item(2326, 7).Value = 1027.3

$ws.Cells.Item(2327, 2).Value = 'M&MFIN'
$ws.Cells.Item(2327, 7).Value = 287.9

$ws.Cells.Item(2328, 2).Value = 'MANORG'
$ws.Cells.Item(2328, 7).Value = 450.8

$ws.Cells.Item(2329, 2).Value = 'MIDHANI'
$ws.Cells.Item(2329, 7).Value = 406.75

$ws.Cells.Item(2330, 2).Value = 'MOMENTUM'
$ws.Cells.Item(2330, 7).Value = 35.56

$ws.Cells.Item(2331, 2).Value = 'MPSLTD'
$ws.Cells.Item(2331, 7).Value = 2119.3

$ws.Cells.Item(2332, 2).Value = 'NAGAFERT'
$ws.Cells.Item(2332, 7).Value = 10.02

$ws.Cells.Item(2333, 2).Value = 'NDTV'
$ws.Cells.Item(2333, 7).Value = 201.72

$ws.Cells.Item(2334, 2).Value = 'NUCLEUS'
$ws.Cells.Item(2334, 7).Value = 1142.2

$ws.Cells.Item(2335, 2).Value = 'NURECA'
$ws.Cells.Item(2335, 7).Value = 261.2

$ws.Cells.Item(2336, 2).Value = 'ORIENTCER'
$ws.Cells.Item(2336, 7).Value = 49.09

$ws.Cells.Item(2337, 2).Value = 'PATINTLOG'
$ws.Cells.Item(2337, 7).Value = 23.07

$ws.Cells.Item(2338, 2).Value = 'PIIND'
$ws.Cells.Item(2338, 7).Value = 4334.3

$ws.Cells.Item(2339, 2).Value = 'PIONEEREMB'
$ws.Cells.Item(2339, 7).Value = 45.24

$ws.Cells.Item(2340, 2).Value = 'PLASTIBLEN'
$ws.Cells.Item(2340, 7).Value = 293.15

$ws.Cells.Item(2341, 2).Value = 'PRSMJOHNSN'
$ws.Cells.Item(2341, 7).Value = 155.63

$ws.Cells.Item(2342, 2).Value = 'RAJMET'
$ws.Cells.Item(2342, 7).Value = 12.25

$ws.Cells.Item(2343, 2).Value = 'RAMANEWS'
$ws.Cells.Item(2343, 7).Value = 19.94

$ws.Cells.Item(2344, 2).Value = 'RATEGAIN'
$ws.Cells.Item(2344, 7).Value = 718

$ws.Cells.Item(2345, 2).Value = 'ROLEXRINGS'
$ws.Cells.Item(2345, 7).Value = 2240.15

$ws.Cells.Item(2346, 2).Value = 'RSWM'
$ws.Cells.Item(2346, 7).Value = 213

$ws.Cells.Item(2347, 2).Value = 'SADBHAV'
$ws.Cells.Item(2347, 7).Value = 28.69

$ws.Cells.Item(2348, 2).Value = 'SAKSOFT'
$ws.Cells.Item(2348, 7).Value = 277.45

$ws.Cells.Item(2349, 1).Value = '14/08/2024'

Write-Output "Added rows 2265-2349"